$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ranking computation was re-run, giving each worker a slightly
# different re_rank score (column F). Rows are re-sorted in descending
# order of that score, while the original "index" (column A) stays
# attached to the worker it was first assigned to (so it is no longer
# sequential after the sort), and the displayed rank (column H) is
# recalculated as the row's new position.
#
# Columns: A=index, B=prolificid, C=prolificid(hash), D=name,
#          E=gender, F=re_rank score, G=race, H=rank

$rows = @(
    @(0, 44, '60b091ed11ccda59e3fc7761', 'Myles',      'male', 13.08405170480421, 'Black or African American'),
    @(1, 3,  '601d69a993d94008fb2b25dc', 'Quinterius', 'male', 8.317598354541474, 'Black or African American'),
    @(2, 30, '60c2341fe95d71ee52c043f0', 'Matthew',    'male', 7.382070019746715, 'White'),
    @(3, 27, '5ff8ad350d084e10f500e48a', 'Drew',       'male', 7.27978505289396,  'White'),
    @(4, 32, '60bf9943e4e04642d4634ecc', 'Jamarii',    'male', 5.430038462157364, 'Black or African American'),
    @(5, 26, '5dd671942b033b5ec8bc97b4', 'Juan',       'male', 5.411470426993446, 'Hispanic'),
    @(6, 22, '60db4fde6193c50664c9c478', 'Edosagbe',   'male', 5.05463621270477,  'Black or African American'),
    @(7, 2,  '60b322994d0b901954690036', 'Corey',      'male', 4.430868679986358, 'White'),
    @(8, 33, '5e2522d6b734b47915f88275', 'Brennan',    'male', 4.412626648038093, 'White'),
    @(9, 49, '6088fc724afd5c008db33e9d', 'Masuf',      'male', 3.327095999247362, 'Asian'),
    @(10, 50, '6097b95056caf5ebb2720002', 'Damian',    'male', 2.26493513038394,  'Black or African American'),
    @(11, 29, '60b83826821417f8e484a207', 'Eli',       'male', 2.202334476874346, 'White')
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $i + 1
}
